# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Note: some Price values are leading-quoted ('1.002, '214.48, ...) so Excel
# stores them as text instead of auto-coercing the single-dot numeric-looking
# strings into numbers (multi-dot values like 25.952.47 are never coerced and
# need no prefix). This mirrors the original inline-string cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.952.47'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.638.67'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -1.12%  '
$ws.Range('D5').Value = '''214.48'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '''0.5099'
$ws.Range('E6').Value = '  +1.58%  '
$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').Value = '''0.2558'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '''0.06357'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '''19.51'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = '''0.07757'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '''4.287'
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = '1.650.83'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').Value = '''0.5437'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = '0.0₅7747'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '''64.22'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '25.976.39'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '''196.27'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').Value = '''4.429'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('D21').Value = '''9.917'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').Value = '''6.047'
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').Value = '''1.002'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').Value = '''1.869'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('D25').Value = '''141.03'
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = '''0.1193'
$ws.Range('E26').Value = '  +5.91%  '
$ws.Range('D27').Value = '''6.837'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').Value = '''15.60'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = '''1.235'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').Value = '''0.04939'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').Value = '''3.245'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').Value = '''3.177'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = '''1.524'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('D35').Value = '''0.8917'
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('D36').Value = '1.148.51'
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('D37').Value = '''2.587'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('D38').Value = '''0.5429'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').Value = '''0.01553'
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').Value = '''1.001'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('D42').Value = '0.0₈127'
$ws.Range('E42').Value = '  +5.51%  '
$ws.Range('D43').Value = '''0.8124'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').Value = '''99.02'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '''5.452'
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('D46').Value = '1.774.78'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = '''0.4525'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '''1.000'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''54.80'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').Value = '''0.05055'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = '''1.002'
$ws.Range('E51').Value = '  -0.40%  '
